{"js": "// Update the PO1 value figure from 43.471.943.900 to 43.377.130.500.\nconst body = context.document.body;\n\nconst searchResults = body.search(\"43.471.943.900\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < searchResults.items.length; i++) {\n  searchResults.items[i].insertText(\"43.377.130.500\", Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Update the PO1 value figure from 43.471.943.900 to 43.377.130.500.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"43.471.943.900\"\n$find.Replacement.Text = \"43.377.130.500\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n\n$wdReplaceAll = 2\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, $wdReplaceAll) | Out-Null\n"}
